$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the helper SUM formula in D1 (column D no longer used)
$ws.Range("D1").ClearContents()

# 2. Reorder the three "# ..." count labels in rows 3-5
#    old: row3="# alt formes", row4="# evolves into", row5="# evolves from"
#    new: row3="# evolves into", row4="# evolves from", row5="# alt formes"
$ws.Range("B3").Value = "# evolves into"
$ws.Range("B4").Value = "# evolves from"
$ws.Range("B5").Value = "# alt formes"

# 3. Update the descriptive "List of ..." rows (66-70), reshuffled:
#    new row66: egg moves list description (text revised: "Move" -> "Move index")
$ws.Range("B66").Value = "List of egg moves Move index (2 bytes) each"

#    new row67: "evolves into" description in B, and the old levelup-moves description moves to C67
$ws.Range("B67").Value = "List of Pokemon it evolves into: Method (1 byte) target (2) method (2) forme target (1) number parameter (1)"
$ws.Range("C67").Value = "List of levelup moves Move (2 bytes) level (1 byte) each"

#    new row68: "evolves into it" description
$ws.Range("B68").Value = "List of evolves into it: source (2) forme (1) Method (1 byte) method (2)  number parameter (1)"

#    new row69: "alt formes" description
$ws.Range("B69").Value = "List of alt formes: forme # (1), method of transformation (1), item or move needed (2)"

# 4. The old row 70 is no longer needed; delete it so the sheet ends at row 69
$ws.Rows(70).Delete()

# 5. Leave the cursor on C1, matching the saved selection state
$null = $ws.Range("C1").Select()
